$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3572712
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 3704997.8
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 11114993.4
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -11115329.4

$ws.Range("H76").Value = 16938.281
$ws.Range("I76").Value = 10239
$ws.Range("J76").Value = 17923.47
$ws.Range("K76").Value = 10239
$ws.Range("L76").Value = 17923.47
$ws.Range("M76").Value = -9924
$ws.Range("N76").Value = -18553.47

$ws.Range("H79").Value = 16938.281
$ws.Range("I79").Value = 10239
$ws.Range("J79").Value = 17923.47
$ws.Range("K79").Value = 10239
$ws.Range("L79").Value = 17923.47
$ws.Range("M79").Value = -9147
$ws.Range("N79").Value = -20107.47

$ws.Range("H86").Value = 3949
$ws.Range("I86").Value = 3684.2856
$ws.Range("K86").Value = 3684.2856
$ws.Range("M86").Value = -2561.2856

$ws.Range("H89").Value = 3949
$ws.Range("I89").Value = 3684.2856
$ws.Range("K89").Value = 18421.428
$ws.Range("M89").Value = -12805.428

$ws.Range("H111").Value = 3007.2
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 3007.2
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 9021.599999999999
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -15155.6

$ws.Range("H137").Value = 1952.5
$ws.Range("I137").Value = 1799.6428
$ws.Range("K137").Value = 5398.928400000001
$ws.Range("M137").Value = -2848.928400000001

$ws.Range("H138").Value = 2796.5532
$ws.Range("J138").Value = 2665.6316
$ws.Range("L138").Value = 7996.8948
$ws.Range("N138").Value = -18276.8948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1600
$ws.Range("I2").Value = 1600
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1600
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1487
$ws.Range("N2").ClearContents()

$ws.Range("H61").Value = 8643.032999999999
$ws.Range("I61").Value = 7451.68
$ws.Range("K61").Value = 7451.68
$ws.Range("M61").Value = -7239.68

$ws.Range("H116").Value = 1600
$ws.Range("I116").Value = 1600
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1600
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 694
$ws.Range("N116").ClearContents()

$ws.Range("H136").Value = 8643.032999999999
$ws.Range("I136").Value = 7451.68
$ws.Range("K136").Value = 22355.04
$ws.Range("M136").Value = -19805.04

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 62753.332
$ws.Range("J2").Value = 69130
$ws.Range("L2").Value = 69130
$ws.Range("N2").Value = -69356

$ws.Range("H3").Value = 1600
$ws.Range("I3").Value = 1600
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1600
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1486
$ws.Range("N3").ClearContents()

$ws.Range("H7").Value = 4997
$ws.Range("I7").Value = 4995.5
$ws.Range("K7").Value = 4995.5
$ws.Range("M7").Value = -4882.5

$ws.Range("H20").Value = 3767.0557
$ws.Range("I20").Value = 3383.8
$ws.Range("J20").Value = 4246.125
$ws.Range("K20").Value = 3383.8
$ws.Range("L20").Value = 4246.125
$ws.Range("M20").Value = -3136.8
$ws.Range("N20").Value = -4740.125

$ws.Range("H86").Value = 1923.9429
$ws.Range("I86").Value = 1683.2413
$ws.Range("K86").Value = 1683.2413
$ws.Range("M86").Value = -560.2412999999999

$ws.Range("H89").Value = 1923.9429
$ws.Range("I89").Value = 1683.2413
$ws.Range("K89").Value = 8416.2065
$ws.Range("M89").Value = -2800.2065

$ws.Range("H99").Value = 3848.923
$ws.Range("I99").Value = 2782.2666
$ws.Range("K99").Value = 2782.2666
$ws.Range("M99").Value = -1284.2666

$ws.Range("H132").Value = 100000
$ws.Range("J132").Value = 100000
$ws.Range("L132").Value = 100000
$ws.Range("N132").Value = -110120

$ws.Range("H134").Value = 4418.9375
$ws.Range("I134").Value = 4216.567
$ws.Range("J134").Value = 7454.5
$ws.Range("K134").Value = 12649.701
$ws.Range("L134").Value = 22363.5
$ws.Range("M134").Value = -10114.701
$ws.Range("N134").Value = -27433.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3540.0417
$ws.Range("I16").Value = 1893.1875
$ws.Range("K16").Value = 1893.1875
$ws.Range("M16").Value = -1606.1875

$ws.Range("H113").Value = 3540.0417
$ws.Range("I113").Value = 1893.1875
$ws.Range("K113").Value = 1893.1875
$ws.Range("M113").Value = 276.8125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2684.8
$ws.Range("I68").Value = 3321.4285
$ws.Range("J68").Value = 1199.3334
$ws.Range("K68").Value = 9964.2855
$ws.Range("L68").Value = 3598.0002
$ws.Range("M68").Value = -9153.2855
$ws.Range("N68").Value = -5220.0002

$ws.Range("H71").Value = 2684.8
$ws.Range("I71").Value = 3321.4285
$ws.Range("J71").Value = 1199.3334
$ws.Range("K71").Value = 29892.8565
$ws.Range("L71").Value = 10794.0006
$ws.Range("M71").Value = -25836.8565
$ws.Range("N71").Value = -18906.0006

$ws.Range("H86").Value = 608.3333
$ws.Range("I86").Value = 512.5
$ws.Range("K86").Value = 1537.5
$ws.Range("M86").Value = -351.5

$ws.Range("H89").Value = 608.3333
$ws.Range("I89").Value = 512.5
$ws.Range("K89").Value = 4612.5
$ws.Range("M89").Value = 1315.5

$ws.Range("H128").Value = 1063262.9
$ws.Range("I128").Value = 1063262.9
$ws.Range("K128").Value = 3189788.7
$ws.Range("M128").Value = -3184808.7

$ws.Range("H132").Value = 1931.9524
$ws.Range("I132").Value = 1708.4445
$ws.Range("K132").Value = 15376.0005
$ws.Range("M132").Value = -12846.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12356.833
$ws.Range("I70").Value = 9587.286
$ws.Range("K70").Value = 9587.286
$ws.Range("M70").Value = -9317.286

$ws.Range("H73").Value = 12356.833
$ws.Range("I73").Value = 9587.286
$ws.Range("K73").Value = 9587.286
$ws.Range("M73").Value = -8651.286

$ws.Range("H96").Value = 30051.4
$ws.Range("J96").Value = 30051.4
$ws.Range("L96").Value = 30051.4
$ws.Range("N96").Value = -35543.4

$ws.Range("H122").Value = 5454.905
$ws.Range("I122").Value = 5883.727
$ws.Range("J122").Value = 4983.2
$ws.Range("K122").Value = 17651.181
$ws.Range("L122").Value = 14949.6
$ws.Range("M122").Value = -15201.181
$ws.Range("N122").Value = -19849.6

$ws.Range("H132").Value = 1218.3334
$ws.Range("I132").Value = 1218.3334
$ws.Range("K132").Value = 3655.0002
$ws.Range("M132").Value = -1125.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2616.4443
$ws.Range("I7").Value = 2781
$ws.Range("J7").Value = 1300
$ws.Range("K7").Value = 2781
$ws.Range("L7").Value = 1300
$ws.Range("M7").Value = -2669
$ws.Range("N7").Value = -1524

$ws.Range("H16").Value = 1596.409
$ws.Range("I16").Value = 1490.0555
$ws.Range("K16").Value = 1490.0555
$ws.Range("M16").Value = -1320.0555

$ws.Range("H22").Value = 1711.6666
$ws.Range("I22").Value = 1554
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 1554
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -1259
$ws.Range("N22").Value = -3090

$ws.Range("H27").Value = 1711.6666
$ws.Range("I27").Value = 1554
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 1554
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -1447
$ws.Range("N27").Value = -2714

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H93").Value = 12184.9
$ws.Range("I93").Value = 2544.2222
$ws.Range("J93").Value = 20072.727
$ws.Range("K93").Value = 2544.2222
$ws.Range("L93").Value = 20072.727
$ws.Range("M93").Value = -1296.2222
$ws.Range("N93").Value = -22568.727

$ws.Range("H126").Value = 2616.4443
$ws.Range("I126").Value = 2781
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 8343
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -5873
$ws.Range("N126").Value = -8840

$ws.Range("H136").Value = 3644.5217
$ws.Range("I136").Value = 3033.3901
$ws.Range("K136").Value = 9100.1703
$ws.Range("M136").Value = -6550.1703

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 11500
$ws.Range("J34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("N34").Value = -10406

$ws.Range("H100").Value = 619.3333
$ws.Range("J100").Value = 347.66666
$ws.Range("L100").Value = 695.33332
$ws.Range("N100").Value = -1777.33332

$ws.Range("H132").Value = 4129.2705
$ws.Range("I132").Value = 4850.625
$ws.Range("J132").Value = 2797.5386
$ws.Range("K132").Value = 14551.875
$ws.Range("L132").Value = 8392.6158
$ws.Range("M132").Value = -12021.875
$ws.Range("N132").Value = -13452.6158
